$wb = $excel.ActiveWorkbook

# --- Sheet 2: "Forecast Scopes Series" -> add row 6 (Prelevato / A_PuntoDiPrelievo block) ---
$wsSeries = $wb.Worksheets.Item("Forecast Scopes Series")
$wsSeries.Select()

# The new row re-uses the plain/unstyled look of hiddenSheet -> pull its format
# (no fill/font override, "Normal" look) onto the new row before filling values.
$wsHidden = $wb.Worksheets.Item("hiddenSheet")
$wsHidden.Range("A1").Copy()
$wsSeries.Range("A6:H6").PasteSpecial(-4122)

$seriesVals = @("CREATE/MODIFY","Prelevato","A_PuntoDiPrelievo","A_IMP_TOT_PRELEVATO","IMP_TOT_PRELEVATO_Auto","IMP_TOT_PRELEVATO_Man","IMP_TOT_PRELEVATO_Flag","IMP_TOT_PRELEVATO_Mix")
for ($i = 0; $i -lt 8; $i++) {
    $wsSeries.Cells.Item(6, $i + 1).Value = $seriesVals[$i]
}
$wsSeries.Range("A7:XFD8").Select()

# --- Sheet 3: "Forecast Scopes Quantiles" -> add rows 5 and 6 ---
$wsQuant = $wb.Worksheets.Item("Forecast Scopes Quantiles")
$wsQuant.Select()

$wsQuant.Cells.Item(5, 1).Value = "CREATE/MODIFY"
$wsQuant.Cells.Item(5, 2).Value = "Prelevato"
$wsQuant.Cells.Item(5, 3).Value = "A_PuntoDiPrelievo"
$wsQuant.Cells.Item(5, 4).Value = "1;2;5;95;98;99"
$wsQuant.Cells.Item(5, 1).Style = "Normal"
$wsQuant.Cells.Item(2, 3).Copy()
$wsQuant.Cells.Item(5, 3).PasteSpecial(-4122)
$wsQuant.Cells.Item(2, 4).Copy()
$wsQuant.Cells.Item(5, 4).PasteSpecial(-4122)

$wsQuant.Cells.Item(6, 1).Value = "CREATE/MODIFY"
$wsQuant.Cells.Item(6, 2).Value = "PrelevatoAggr"
$wsQuant.Cells.Item(6, 3).Value = "A_PuntoDiPrelievo"
$wsQuant.Cells.Item(6, 4).Value = "1;2;5;95;98;99"
$wsQuant.Cells.Item(6, 1).Style = "Normal"
$wsQuant.Cells.Item(2, 3).Copy()
$wsQuant.Cells.Item(6, 3).PasteSpecial(-4122)
$wsQuant.Cells.Item(2, 4).Copy()
$wsQuant.Cells.Item(6, 4).PasteSpecial(-4122)
$wsQuant.Range("A7:XFD7").Select()

# --- Sheet 7 (tab): "Forecast Scopes Labels" -> move selection only ---
$wsLabels = $wb.Worksheets.Item("Forecast Scopes Labels")
$wsLabels.Select()
$wsLabels.Range("A5:XFD5").Select()

# --- Sheet 1: "Forecast Scopes" -> move selection and make it the active tab ---
$wsScopes = $wb.Worksheets.Item("Forecast Scopes")
$wsScopes.Select()
$wsScopes.Range("A5:XFD5").Select()
